# Update the "last updated" timestamp in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 20:05"

# --- Simple numeric refreshes (country stays in the same row) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1675880
$ws.Range("C4").Value = 9052
$ws.Range("D4").Value = 451124
$ws.Range("E4").Value = 1125753
$ws.Range("G4").Value = 320
$ws.Range("H4").Value = 99003

# Row 5: Brasil
$ws.Range("B5").Value = 352744
$ws.Range("C5").Value = 5346
$ws.Range("E5").Value = 187866
$ws.Range("G5").Value = 278
$ws.Range("H5").Value = 22291

# Row 7: España
$ws.Range("B7").Value = 282852
$ws.Range("C7").Value = 482
$ws.Range("E7").Value = 57142
$ws.Range("G7").Value = 74
$ws.Range("H7").Value = 28752

# Row 13: India
$ws.Range("B13").Value = 138223
$ws.Range("C13").Value = 6800
$ws.Range("D13").Value = 57643
$ws.Range("E13").Value = 76556
$ws.Range("G13").Value = 156
$ws.Range("H13").Value = 4024

# Row 16: Canada
$ws.Range("B16").Value = 84655
$ws.Range("C16").Value = 1034
$ws.Range("D16").Value = 43927
$ws.Range("E16").Value = 34304
$ws.Range("G16").Value = 69
$ws.Range("H16").Value = 6424

# Row 30: Suiza
$ws.Range("D30").Value = 28100
$ws.Range("E30").Value = 730

# Row 41: Israel
$ws.Range("B41").Value = 16717
$ws.Range("C41").Value = 5
$ws.Range("D41").Value = 14153
$ws.Range("E41").Value = 2285

# Row 53: Barein
$ws.Range("B53").Value = 9138
$ws.Range("C53").Value = 336
$ws.Range("D53").Value = 4587
$ws.Range("E53").Value = 4538

# Row 146: Montenegro
$ws.Range("D146").Value = 315
$ws.Range("E146").Value = 0

# Row 153: Yemen
$ws.Range("B153").Value = 222
$ws.Range("C153").Value = 10
$ws.Range("D153").Value = 10
$ws.Range("E153").Value = 170
$ws.Range("G153").Value = 3
$ws.Range("H153").Value = 42

# Row 179: Angola
$ws.Range("B179").Value = 69
$ws.Range("C179").Value = 8
$ws.Range("E179").Value = 47

# --- Re-sorted rows: Sudafrica overtakes Indonesia (rows 34/35) ---
# Row 34 becomes Sudafrica with its freshly updated numbers
$ws.Range("A34").Value = "Sudafrica"
$ws.Range("B34").Value = 22583
$ws.Range("C34").Value = 1240
$ws.Range("D34").Value = 10104
$ws.Range("E34").Value = 12050
$ws.Range("G34").Value = 22
$ws.Range("H34").Value = 429

# Row 35 becomes Indonesia, keeping its previous (unrevised) numbers
$ws.Range("A35").Value = "Indonesia"
$ws.Range("B35").Value = 22271
$ws.Range("C35").Value = 526
$ws.Range("D35").Value = 5402
$ws.Range("E35").Value = 15497
$ws.Range("G35").Value = 21
$ws.Range("H35").Value = 1372

# --- Re-sorted rows: Republica del Chad overtakes San Marino & Sudan del Sur (rows 126/127/128) ---
# Row 126 becomes Republica del Chad with its freshly updated numbers
$ws.Range("A126").Value = "Republica del Chad"
$ws.Range("B126").Value = 675
$ws.Range("C126").Value = 27
$ws.Range("D126").Value = 215
$ws.Range("E126").Value = 400
$ws.Range("H126").Value = 60

# Row 127 becomes San Marino, keeping its previous (unrevised) numbers
$ws.Range("A127").Value = "San Marino"
$ws.Range("B127").Value = 665
$ws.Range("D127").Value = 266
$ws.Range("E127").Value = 357
$ws.Range("H127").Value = 42

# Row 128 becomes Sudan del Sur, keeping its previous (unrevised) numbers
$ws.Range("A128").Value = "Sudan del Sur"
$ws.Range("B128").Value = 655
$ws.Range("D128").Value = 6
$ws.Range("E128").Value = 641
$ws.Range("H128").Value = 8
